$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 31 de Mayo de 2020 a las 16:35"

# Estados Unidos (row 4): refresh case counts
$ws.Range("B4").Value = 1819792
$ws.Range("C4").Value = 2972
$ws.Range("D4").Value = 535379
$ws.Range("E4").Value = 1178779

# Reino Unido (row 8): refresh case counts
$ws.Range("B8").Value = 274762
$ws.Range("C8").Value = 1936
$ws.Range("G8").Value = 113
$ws.Range("H8").Value = 38489

# India (row 11): refresh case counts
$ws.Range("B11").Value = 186186
$ws.Range("C11").Value = 4359
$ws.Range("D11").Value = 88769
$ws.Range("E11").Value = 92149
$ws.Range("G11").Value = 83
$ws.Range("H11").Value = 5268

# Rumania (row 41): refresh case counts
$ws.Range("E41").Value = 4735
$ws.Range("G41").Value = 7
$ws.Range("H41").Value = 1266

# Israel and Republica Dominicana swap ranking positions (rows 43/44):
# Republica Dominicana moves up to row 43 with fresh numbers, Israel drops
# to row 44 keeping its previous numbers.
$ws.Range("A43").Value = "Republica Dominicana"
$ws.Range("B43").Value = 17285
$ws.Range("C43").Value = 377
$ws.Range("D43").Value = 9557
$ws.Range("E43").Value = 7226
$ws.Range("G43").Value = 4
$ws.Range("H43").Value = 502

$ws.Range("A44").Value = "Israel"
$ws.Range("B44").Value = 17024
$ws.Range("C44").Value = 12
$ws.Range("D44").Value = 14812
$ws.Range("E44").Value = 1928
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 284

# Kazajistan (row 55): refresh case counts
$ws.Range("D55").Value = 5404
$ws.Range("E55").Value = 5414

# Hungria and Tayikistan swap ranking positions (rows 75/76):
# Tayikistan moves up to row 75 with fresh numbers, Hungria drops to row 76
# keeping its previous numbers.
$ws.Range("A75").Value = "Tayikistan"
$ws.Range("B75").Value = 3930
$ws.Range("C75").Value = 123
$ws.Range("D75").Value = 2004
$ws.Range("E75").Value = 1879
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 47

$ws.Range("A76").Value = "Hungria"
$ws.Range("B76").Value = 3876
$ws.Range("C76").Value = 9
$ws.Range("D76").Value = 2147
$ws.Range("E76").Value = 1203
$ws.Range("G76").Value = 2
$ws.Range("H76").Value = 526

# Kenia (row 92): refresh case counts
$ws.Range("D92").Value = 478
$ws.Range("E92").Value = 1420
$ws.Range("G92").Value = 1
$ws.Range("H92").Value = 64
